$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete rows 2 and 3 (TP-9 and TP-8 entries), keeping row 4 (TP-7) which
# shifts up to become the new row 2.
$ws.Rows("2:3").Delete()
